$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.448.98'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '1.916.48'
$ws.Range('E3').Value = '  +1.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4697'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2872'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06821'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '110.13'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.39'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07732'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').Value = '1.888.58'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.307'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6580'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '295.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('D17').Value = '30.450.69'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007625'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '2.144.96'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9994'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.245'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.196'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.383'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').Value = '  +5.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.104'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.91%  '
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.364'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.181'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.992'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05029'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7366'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.153'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02077'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.741'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.679'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '109.67'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8717'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.863'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4262'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '51.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +19.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.187'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.256'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1219'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2463'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.09%  '
